$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column D
$ws.Range("D1").Value = "Age_sex"

# Fill D2:D9 with "Sex Age" combination (matches B<row> & " " & A<row>)
for ($r = 2; $r -le 9; $r++) {
    $age = $ws.Cells.Item($r, 1).Text
    $sex = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 4).Value = "$sex $age"
}

# Update selection to G10, matching the saved sheet view state
$ws.Range("G10").Select()
